# Update "想去人数" (F column) figures for several rows across the
# "展览" (Exhibition), "本地生活" (Local Life) and "全部类型" (All Types) sheets,
# matching the refreshed data pulled from the source site.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsLocalLife  = $wb.Worksheets.Item("本地生活")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet updates
$wsExhibition.Range("F6").Value  = 1139
$wsExhibition.Range("F7").Value  = 2253
$wsExhibition.Range("F8").Value  = 2177
$wsExhibition.Range("F9").Value  = 1140
$wsExhibition.Range("F12").Value = 1703
$wsExhibition.Range("F13").Value = 414
$wsExhibition.Range("F17").Value = 253
$wsExhibition.Range("F18").Value = 1611
$wsExhibition.Range("F19").Value = 272
$wsExhibition.Range("F20").Value = 1332
$wsExhibition.Range("F21").Value = 752
$wsExhibition.Range("F22").Value = 279
$wsExhibition.Range("F23").Value = 634
$wsExhibition.Range("F24").Value = 12409
$wsExhibition.Range("F25").Value = 12455
$wsExhibition.Range("F27").Value = 714
$wsExhibition.Range("F31").Value = 415
$wsExhibition.Range("F32").Value = 1939
$wsExhibition.Range("F35").Value = 214
$wsExhibition.Range("F36").Value = 628

# 本地生活 (Local Life) sheet updates
$wsLocalLife.Range("F3").Value = 112

# 全部类型 (All Types) sheet updates
$wsAllTypes.Range("F7").Value  = 1139
$wsAllTypes.Range("F8").Value  = 2253
$wsAllTypes.Range("F9").Value  = 2177
$wsAllTypes.Range("F10").Value = 1140
$wsAllTypes.Range("F12").Value = 112
$wsAllTypes.Range("F14").Value = 1703
$wsAllTypes.Range("F15").Value = 414
$wsAllTypes.Range("F22").Value = 253
$wsAllTypes.Range("F23").Value = 1611
$wsAllTypes.Range("F24").Value = 272
$wsAllTypes.Range("F25").Value = 1332
$wsAllTypes.Range("F26").Value = 752
$wsAllTypes.Range("F27").Value = 279
$wsAllTypes.Range("F29").Value = 634
$wsAllTypes.Range("F30").Value = 12409
$wsAllTypes.Range("F31").Value = 12455
$wsAllTypes.Range("F33").Value = 714
$wsAllTypes.Range("F37").Value = 415
$wsAllTypes.Range("F40").Value = 1939
$wsAllTypes.Range("F45").Value = 214
$wsAllTypes.Range("F46").Value = 628
